# Generate Report for Archive
# The localization status for the single tracked file moved on from
# "Ready for handoff" to "In Translation" - update every sheet that
# surfaces the Status column (Overview's per-locale summary columns and
# each locale sheet's own Status column), then let the Status column
# re-fit its new (shorter) contents.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Columns previously sized for "Ready for handoff" can now shrink to fit
# the shorter "In Translation" text.
$overview.Columns("E:F").ColumnWidth = 12.576851254417766
$zhcn.Columns("C").ColumnWidth = 12.576851254417766
$dede.Columns("C").ColumnWidth = 12.576851254417766
